$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.768.08'
$ws.Range("E2").Value = '  +2.03%  '
$ws.Range("D3").Value = '3.033.14'
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.83%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.026.23'
$ws.Range("E8").Value = '  +0.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.518'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.63'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +12.44%  '
$ws.Range("E11").Value = '  +2.06%  '
$ws.Range("E12").Value = '  +1.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000234'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.59'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.127'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.86%  '
$ws.Range("D16").Value = '3.534.50'
$ws.Range("E16").Value = '  +1.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").Value = '62.789.14'
$ws.Range("E18").Value = '  +2.03%  '
$ws.Range("D19").Value = '3.029.05'
$ws.Range("E19").Value = '  +0.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '452.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.697'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.71%  '
$ws.Range("E26").Value = '  +2.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.60%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.72%  '
$ws.Range("E30").Value = '  +0.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.51%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("E34").Value = '  +0.80%  '
$ws.Range("D35").Value = '0.0₃0864'
$ws.Range("E35").Value = '  +2.98%  '
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.92'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.18%  '
$ws.Range("E38").Value = '  +8.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.40'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("E41").Value = '  +2.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.10'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.301'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +10.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '394.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0359'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.45%  '
$ws.Range("D47").Value = '2.725.48'
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.73'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.49%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("E50").Value = '  +4.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.40'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.80%  '
